$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Age of The Gods Prince of Olympus Free | Slot
#    Review").
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = '<w:p ' + $wns + '>' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Discover the world of the Greek hero Hercules and try your luck with exciting bonus rounds and Progressive Jackpots in Age of The Gods Prince of Olympus. Play for free now!</w:t></w:r>' +
           '</w:p>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated bold title paragraph
#    ("Play Age of The Gods Prince of Olympus Free | Slot Review") and
#    replace the text of the italic paragraph that follows it with the new
#    image-generation prompt.
# ---------------------------------------------------------------------------
$dupTarget = "Play Age of The Gods Prince of Olympus Free | Slot Review"
$dupPara = $null
for ($i = 3; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs($i)
  $txt = $para.Range.Text.TrimEnd("`r")
  if ($txt -eq $dupTarget) {
    $dupPara = $para
  }
}
$dupPara.Range.Delete()

$descTarget = "Discover the world of the Greek hero Hercules and try your luck with exciting bonus rounds and Progressive Jackpots in Age of The Gods Prince of Olympus. Play for free now!"
$promptText = 'Prompt: Create a feature image for "Age of The Gods Prince of Olympus" that captures the essence of the game - the powerful son of Zeus - Hercules, and the fun and excitement that this slot game brings. Specifications: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses - The image should be eye-catching and playful - The image should include the game title "Age of The Gods Prince of Olympus" Keep in mind the theme of the game, which revolves around Greek mythology and Hercules. The Maya warrior with glasses adds an element of surprise and fun to the image, resulting in a unique and captivating piece of art. The image should aim to draw players in and encourage them to try out the game.'

$descPara = $null
for ($i = 3; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs($i)
  $txt = $para.Range.Text.TrimEnd("`r")
  if ($txt -eq $descTarget) {
    $descPara = $para
  }
}
$descRange = $descPara.Range
$descRange.MoveEnd(1, -1) | Out-Null
$descRange.Text = $promptText

Write-Output "done"
